# Updated cryptos list on Fri Aug 25 18:35:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.018.06'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '1.649.50'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  -0.39%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '218.20'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').Value = '  -0.32%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.2619'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  -1.83%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.06290'
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '20.26'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  -4.01%  '
$ws.Range('E11').Value = '  +0.13%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '4.469'
$r.Style = 'Normal'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '1.600.95'
$ws.Range('E13').Value = '  -3.29%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '0.5451'
$r.Style = 'Normal'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '0.0₅8083'
$ws.Range('E15').Value = '  -1.78%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '64.99'
$r.Style = 'Normal'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '26.036.03'
$ws.Range('E17').Value = '  -0.61%  '
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '1.003'
$r.Style = 'Normal'
$ws.Range('E18').Value = '  -0.30%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '4.556'
$r.Style = 'Normal'
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '192.73'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  +0.05%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '10.04'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('E23').Value = '  -0.47%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '139.34'
$r.Style = 'Normal'
$ws.Range('E24').Value = '  +1.62%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '0.1242'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  +0.27%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '7.267'
$r.Style = 'Normal'
$ws.Range('E26').Value = '  +0.39%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '16.24'
$r.Style = 'Normal'
$ws.Range('E27').Value = '  +0.74%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '1.416'
$r.Style = 'Normal'
$ws.Range('E28').Value = '  +0.20%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '0.05941'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  -1.57%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '1.275'
$r.Style = 'Normal'
$ws.Range('E30').Value = '  -0.57%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '3.489'
$r.Style = 'Normal'
$ws.Range('E31').Value = '  -2.19%  '
$ws.Range('E32').Value = '  -2.75%  '
$ws.Range('E33').Value = '  -6.96%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '2.412'
$r.Style = 'Normal'
$ws.Range('E34').Value = '  +0.02%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '0.9423'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  -3.94%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '2.745'
$r.Style = 'Normal'
$ws.Range('E36').Value = '  -0.99%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.5655'
$r.Style = 'Normal'
$ws.Range('E37').Value = '  -4.92%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.01607'
$r.Style = 'Normal'
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('E39').Value = '  -1.69%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.8483'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('E41').Value = '  -0.29%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '101.00'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  +1.33%  '
$ws.Range('D43').Value = '1.006.51'
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('D44').Value = '1.787.87'
$ws.Range('E44').Value = '  -0.39%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '56.74'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -5.22%  '
$ws.Range('E47').Value = '  +0.05%  '
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '1.481'
$r.Style = 'Normal'
$ws.Range('E48').Value = '  +1.02%  '
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '0.4284'
$r.Style = 'Normal'
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '7.851'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.05149'
$r.Style = 'Normal'
$ws.Range('E51').Value = '  -0.57%  '
